# Update cryptos list worksheet with latest prices/volumes (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to stay plain text (matches source data which always
    # stores Price column as text, even when it looks numeric) and avoid
    # leaving a lingering custom number format on the cell afterwards.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "24.825.48"
$ws.Range("E2").Value = "  +0.95%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.707.15"
$ws.Range("E3").Value = "  +1.10%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.27%  "

# Row 5 - BNB
Set-TextValue "D5" "315.00"
$ws.Range("E5").Value = "  +0.65%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.15%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4001"
$ws.Range("E7").Value = "  +2.69%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.4039"
$ws.Range("E8").Value = "  +0.43%  "

# Row 9 - now BinanceUSD (was Polygon)
$ws.Range("B9").Value = "BinanceUSD"
$ws.Range("C9").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D9" "1.002"
$ws.Range("E9").Value = "  -0.32%  "

# Row 10 - now Polygon (was BinanceUSD)
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D10" "1.473"
$ws.Range("E10").Value = "  -1.53%  "

# Row 11 - OKB
$ws.Range("E11").Value = "  +1.98%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.08812"
$ws.Range("E12").Value = "  +0.80%  "

# Row 13 - Solana
Set-TextValue "D13" "26.34"
$ws.Range("E13").Value = "  +6.41%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.516"
$ws.Range("E14").Value = "  -0.93%  "

# Row 15 - Chainlink
Set-TextValue "D15" "8.000"
$ws.Range("E15").Value = "  +0.59%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.00001340"
$ws.Range("E16").Value = "  -0.63%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "1.751.68"
$ws.Range("E17").Value = "  +4.17%  "

# Row 18 - Litecoin
Set-TextValue "D18" "95.54"
$ws.Range("E18").Value = "  -2.78%  "

# Row 19 - TRON
Set-TextValue "D19" "0.07167"
$ws.Range("E19").Value = "  +1.05%  "

# Row 20 - Avalanche
Set-TextValue "D20" "20.91"
$ws.Range("E20").Value = "  +5.58%  "

# Row 21 - Uniswap
Set-TextValue "D21" "7.287"
$ws.Range("E21").Value = "  +0.23%  "

# Row 22 - Dai
Set-TextValue "D22" "1.001"
$ws.Range("E22").Value = "  -0.18%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  +1.37%  "

# Row 24 - WrappedBTC
Set-TextValue "D24" "24.808.35"
$ws.Range("E24").Value = "  +0.87%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.352"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26 - LidoDAOToken
Set-TextValue "D26" "2.895"
$ws.Range("E26").Value = "  -3.49%  "

# Row 27 - HuobiToken
Set-TextValue "D27" "6.389"
$ws.Range("E27").Value = "  +22.62%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +1.70%  "

# Row 29 - Monero
Set-TextValue "D29" "161.69"
$ws.Range("E29").Value = "  +0.20%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "143.79"
$ws.Range("E30").Value = "  +5.55%  "

# Row 31 - Filecoin
Set-TextValue "D31" "8.212"
$ws.Range("E31").Value = "  -4.21%  "

# Row 32 - WEMIXTOKEN
$ws.Range("E32").Value = "  +14.49%  "

# Row 33 - WrappedliquidstakedEther2.0
Set-TextValue "D33" "1.919.26"
$ws.Range("E33").Value = "  +2.63%  "

# Row 34 - now Hedera (was VeChain)
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.08658"
$ws.Range("E34").Value = "  -1.07%  "

# Row 35 - now VeChain (was Hedera)
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D35" "0.03199"
$ws.Range("E35").Value = "  +10.03%  "

# Row 36 - InternetComputer(DFINITY)
Set-TextValue "D36" "7.298"
$ws.Range("E36").Value = "  -0.92%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "1.030"
$ws.Range("E37").Value = "  -0.68%  "

# Row 38 - Algorand
Set-TextValue "D38" "0.2856"
$ws.Range("E38").Value = "  +4.88%  "

# Row 39 - TheSandbox
Set-TextValue "D39" "0.8408"
$ws.Range("E39").Value = "  +7.77%  "

# Row 40 - Stellar
Set-TextValue "D40" "0.09446"
$ws.Range("E40").Value = "  +3.67%  "

# Row 41 - FraxShare
Set-TextValue "D41" "10.74"
$ws.Range("E41").Value = "  -0.18%  "

# Row 42 - Aptos
Set-TextValue "D42" "14.24"
$ws.Range("E42").Value = "  +0.40%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "1.480"
$ws.Range("E43").Value = "  +1.72%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "17.43"
$ws.Range("E44").Value = "  +4.40%  "

# Row 45 - NEARProtocol
Set-TextValue "D45" "2.726"
$ws.Range("E45").Value = "  +5.71%  "

# Row 46 - Decentraland
Set-TextValue "D46" "0.7420"
$ws.Range("E46").Value = "  +3.29%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "4.220"
$ws.Range("E47").Value = "  +0.57%  "

# Row 48 - Flow
Set-TextValue "D48" "1.370"
$ws.Range("E48").Value = "  +2.88%  "

# Row 49 - Frax
$ws.Range("E49").Value = "  -0.17%  "

# Row 50 - Quant
Set-TextValue "D50" "140.72"
$ws.Range("E50").Value = "  +2.25%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.08383"
$ws.Range("E51").Value = "  +5.27%  "
